$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.441.28'
$ws.Range('E2').Value = '  +1.41%  '
$ws.Range('D3').Value = '1.681.64'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9989'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5315'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9993'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2669'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06461'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.33'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07806'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.26%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.512'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.97%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.675.40'
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5628'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.42%  '
$ws.Range('D15').Value = '0.0₅8465'
$ws.Range('E15').Value = '  +6.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = '26.431.56'
$ws.Range('E17').Value = '  +1.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9997'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.850'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '195.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.41'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.403'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9996'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1268'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.515'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.28%  '
$ws.Range('E27').Value = '  +4.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.438'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06210'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.277'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.556'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.470'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.709'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.022'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.790'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.403'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5757'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01647'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.955'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8695'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.96%  '
$ws.Range('D41').Value = '1.058.22'
$ws.Range('E41').Value = '  -1.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9997'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.35'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').Value = '1.829.18'
$ws.Range('E44').Value = '  +1.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.28'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.69%  '
$ws.Range('E46').Value = '  +2.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.152'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.002'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05197'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.061'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.09964'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.45%  '
